# ---------------------------------------------------------------------------
# "new code by shambhu"
#
# Adds a new "DomesticFTA" worksheet after the existing "Login" sheet and
# populates it with a small domestic-funds-transfer test-data table, mirroring
# the author's commit. Login's own data is left untouched (only its
# "tabSelected" view flag moves to the new, now-active sheet, which Excel
# does automatically when a new sheet becomes ActiveSheet).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- add the new worksheet as the last tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "DomesticFTA"

# --- header row + first data row, written in the original authoring order -
# (this keeps shared-string ids allocated in the same order as the source
# workbook, B1 first, then C1, D1, ... )
$new.Range("B1").Value2 = "ReciverBankName"
$new.Range("C1").Value2 = "ReceiverName"
$new.Range("D1").Value2 = "ReceAccNum"
$new.Range("E1").Value2 = "SwiftMsg"
$new.Range("G1").Value2 = "TransferType"
$new.Range("H1").Value2 = "DOT"
$new.Range("I1").Value2 = "TransferDesc"
$new.Range("B2").Value2 = "SBI"
$new.Range("E2").Value2 = "MT103"
$new.Range("G2").Value2 = "Domestic Transfer"
$new.Range("A1").Value2 = "DataBinding"
$new.Range("A2").Value2 = "Data001"
$new.Range("C2").Value2 = "John"
$new.Range("F1").Value2 = "Amount"
$new.Range("A3").Value2 = "Data002"
$new.Range("A4").Value2 = "Data003"
$new.Range("A5").Value2 = "Data004"
$new.Range("A6").Value2 = "Data005"
$new.Range("B3").Value2 = "HDFC"
$new.Range("B4").Value2 = "RBS"
$new.Range("B5").Value2 = "ICICI"
$new.Range("B6").Value2 = "CITI"
$new.Range("C6").Value2 = "Pitter"
$new.Range("C5").Value2 = "Raj"
$new.Range("C4").Value2 = "Smith"
$new.Range("C3").Value2 = "Jacson"

# --- remaining repeating text columns for rows 3-6 -------------------------
$new.Range("E3").Value2 = "MT103"
$new.Range("E4").Value2 = "MT103"
$new.Range("E5").Value2 = "MT103"
$new.Range("E6").Value2 = "MT103"

$new.Range("G3").Value2 = "Domestic Transfer"
$new.Range("G4").Value2 = "Domestic Transfer"
$new.Range("G5").Value2 = "Domestic Transfer"
$new.Range("G6").Value2 = "Domestic Transfer"

$new.Range("I2").Value2 = "TransferDesc"
$new.Range("I3").Value2 = "TransferDesc"
$new.Range("I4").Value2 = "TransferDesc"
$new.Range("I5").Value2 = "TransferDesc"
$new.Range("I6").Value2 = "TransferDesc"

# --- numeric account-number column (plain numbers, no formatting) ---------
$new.Range("D2").Value2 = 1234556655
$new.Range("D3").Value2 = 1234556656
$new.Range("D4").Value2 = 1234556657
$new.Range("D5").Value2 = 1234556658
$new.Range("D6").Value2 = 1234556659

# --- numeric "days" column ---------------------------------------------
$new.Range("F2").Value2 = 8
$new.Range("F3").Value2 = 10
$new.Range("F4").Value2 = 6
$new.Range("F5").Value2 = 9
$new.Range("F6").Value2 = 2

# --- date-of-transfer column: format H2 as a date, then propagate that
#     exact style (same style index) to H3:H6 via copy/paste-format so the
#     workbook doesn't grow one style record per cell --------------------
$new.Range("H2").Value2 = 117924
$new.Range("H3").Value2 = 118655
$new.Range("H4").Value2 = 120481
$new.Range("H5").Value2 = 118655
$new.Range("H6").Value2 = 117559

$new.Range("H2").NumberFormat = "mm-dd-yy"
$new.Range("H2").Copy() | Out-Null
$new.Range("H3:H6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- column widths, approximating Excel's own "best fit" auto-sizing for
#     the now-populated columns -------------------------------------------
$new.Columns.Item(2).ColumnWidth = 16.59
$new.Columns.Item(3).ColumnWidth = 13.31
$new.Columns.Item(4).ColumnWidth = 13.31
$new.Columns.Item(5).ColumnWidth = 11.74
$new.Columns.Item(7).ColumnWidth = 16.45
$new.Columns.Item(8).ColumnWidth = 12.02

# --- view state: new sheet becomes active tab, selection on K13 -----------
$new.Range("K13").Select() | Out-Null
